# Applies the "Updated cryptos list" data refresh to Sheet1.
# For each changed row, Price (D) and/or Volume(1h) (E) text values are updated
# in place. Price values that look numeric (single decimal point) are written
# with a leading apostrophe so Excel keeps them as text, matching the original
# inline-string cell content instead of converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> @{ D = <new price text>; E = <new volume text> }
$updates = [ordered]@{
    2 = @{ D = '66.187.04'; E = '  -4.08%  ' }
    3 = @{ D = '3.546.71' }
    4 = @{ D = '''0.999'; E = '  -0.03%  ' }
    5 = @{ D = '''576.11'; E = '  -6.38%  ' }
    6 = @{ D = '''186.91'; E = '  -2.61%  ' }
    7 = @{ D = '3.542.50'; E = '  -4.41%  ' }
    8 = @{ D = '''0.610'; E = '  -3.86%  ' }
    9 = @{ E = '  +0.04%  ' }
    10 = @{ E = '  -7.11%  ' }
    11 = @{ E = '  -9.88%  ' }
    12 = @{ D = '''52.49'; E = '  -6.98%  ' }
    13 = @{ E = '  -11.30%  ' }
    14 = @{ D = '''9.71'; E = '  -7.75%  ' }
    15 = @{ D = '4.108.11'; E = '  -4.39%  ' }
    16 = @{ D = '3.546.94'; E = '  -4.34%  ' }
    17 = @{ E = '  -1.04%  ' }
    18 = @{ D = '''18.16'; E = '  -5.72%  ' }
    19 = @{ D = '66.025.99'; E = '  -3.98%  ' }
    20 = @{ D = '''12.06'; E = '  -6.31%  ' }
    21 = @{ E = '  -7.41%  ' }
    22 = @{ D = '''390.85' }
    23 = @{ D = '''4.27'; E = '  -7.25%  ' }
    24 = @{ D = '''85.14'; E = '  -4.36%  ' }
    25 = @{ D = '''10.97'; E = '  +0.59%  ' }
    26 = @{ D = '''2.87'; E = '  -5.18%  ' }
    27 = @{ D = '''12.29'; E = '  -3.80%  ' }
    28 = @{ E = '  -0.09%  ' }
    29 = @{ D = '''3.49'; E = '  -6.65%  ' }
    30 = @{ D = '''8.81'; E = '  -8.60%  ' }
    31 = @{ D = '''30.77'; E = '  -6.65%  ' }
    32 = @{ D = '''7.09' }
    33 = @{ D = '''627.54'; E = '  +0.41%  ' }
    34 = @{ D = '''12.07'; E = '  -4.27%  ' }
    35 = @{ D = '''63.30'; E = '  -3.65%  ' }
    36 = @{ E = '  -7.48%  ' }
    37 = @{ E = '  -7.91%  ' }
    38 = @{ E = '  +0.09%  ' }
    39 = @{ D = '''0.392'; E = '  -4.62%  ' }
    40 = @{ D = '0.0₃0757'; E = '  -6.22%  ' }
    41 = @{ D = '''0.998'; E = '  -0.15%  ' }
    42 = @{ E = '  -7.01%  ' }
    43 = @{ D = '2.962.89'; E = '  +3.34%  ' }
    44 = @{ E = '  -7.62%  ' }
    45 = @{ E = '  -5.31%  ' }
    46 = @{ D = '''0.0403'; E = '  -8.69%  ' }
    47 = @{ E = '  -7.52%  ' }
    48 = @{ E = '  -2.61%  ' }
    49 = @{ D = '''137.91'; E = '  -2.43%  ' }
    50 = @{ D = '''8.37'; E = '  -7.65%  ' }
    51 = @{ E = '  -8.51%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
